$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-07-03 Wednesday" "2024-07-04 Thursday"

Replace-Text "684÷5=136, 4" "319÷6=53, 1"
Replace-Text "423÷6=70, 3" "164÷3=54, 2"
Replace-Text "867÷7=123, 6" "315÷2=157, 1"
Replace-Text "402÷4=100, 2" "436÷7=62, 2"
Replace-Text "707÷5=141, 2" "545÷2=272, 1"

Replace-Text "200÷4=50, 0" "457÷9=50, 7"
Replace-Text "754÷6=125, 4" "185÷6=30, 5"
Replace-Text "653÷7=93, 2" "605÷5=121, 0"
Replace-Text "726÷9=80, 6" "687÷8=85, 7"
Replace-Text "985÷5=197, 0" "167÷7=23, 6"

Replace-Text "721÷3=240, 1" "779÷5=155, 4"
Replace-Text "978÷3=326, 0" "694÷6=115, 4"
Replace-Text "549÷4=137, 1" "544÷7=77, 5"
Replace-Text "218÷9=24, 2" "765÷6=127, 3"
Replace-Text "755÷8=94, 3" "556÷7=79, 3"

Replace-Text "856÷7=122, 2" "302÷6=50, 2"
Replace-Text "372÷3=124, 0" "265÷9=29, 4"
Replace-Text "989÷5=197, 4" "410÷2=205, 0"
Replace-Text "899÷8=112, 3" "506÷5=101, 1"
Replace-Text "915÷7=130, 5" "810÷2=405, 0"

Replace-Text "389÷6=64, 5" "134÷6=22, 2"
Replace-Text "376÷5=75, 1" "178÷9=19, 7"
Replace-Text "860÷3=286, 2" "397÷7=56, 5"
Replace-Text "896÷9=99, 5" "982÷2=491, 0"
Replace-Text "594÷9=66, 0" "967÷8=120, 7"
